$d = $word.ActiveDocument
$d.TrackRevisions = $false

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p3 = @'
<w:p><w:r><w:t xml:space="preserve">If you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>don’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> know the syntax</w:t></w:r><w:r><w:t xml:space="preserve"> of a function in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>R</w:t></w:r><w:r><w:t>studio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, use help</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Either by </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>typing</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> ?</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>function_of_interest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">() </w:t></w:r><w:r><w:t>in the console</w:t></w:r><w:r><w:t xml:space="preserve">. Alternatively, </w:t></w:r><w:r><w:t>you can use</w:t></w:r><w:r><w:t xml:space="preserve"> the Help tab in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rstudio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. You can also google the function. In gener</w:t></w:r><w:r><w:t>al, and I am sure you have heard</w:t></w:r><w:r><w:t xml:space="preserve"> it a lot</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> google is your best friend while working with R. If you run into an unknown</w:t></w:r><w:r><w:t xml:space="preserve"> error: google;</w:t></w:r><w:r><w:t xml:space="preserve"> if you know what you want to do but </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>don’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> know how to do it: google. </w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($pkgHeader + $p3 + $pkgFooter)

$p4 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Import muskox dataset (excel file) into </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rstudio</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($pkgHeader + $p4 + $pkgFooter)

$p5 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Possible function:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>read.xlsx(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>)</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($pkgHeader + $p5 + $pkgFooter)

$p10 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>N</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>ote:</w:t></w:r><w:r><w:t xml:space="preserve"> the </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>ggplot2</w:t></w:r><w:r><w:t xml:space="preserve"> package is widely used and I highly recommend learning it. It </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">is </w:t></w:r><w:r><w:t>mainly used</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> in most built-</w:t></w:r><w:r><w:t xml:space="preserve">in plot functions </w:t></w:r><w:r><w:t>of</w:t></w:r><w:r><w:t xml:space="preserve"> other packages and it makes the base function </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>plot</w:t></w:r><w:r><w:t xml:space="preserve"> almost entirely irrelevant once you learn it. </w:t></w:r></w:p>
'@
$d.Paragraphs(10).Range.InsertXML($pkgHeader + $p10 + $pkgFooter)

$p15 = @'
<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Possible functions:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>subset(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>, for(…)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$d.Paragraphs(15).Range.InsertXML($pkgHeader + $p15 + $pkgFooter)

$p16 = @'
<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Note:</w:t></w:r><w:r><w:t xml:space="preserve"> This is maybe the first </w:t></w:r><w:r><w:t xml:space="preserve">tricky </w:t></w:r><w:r><w:t xml:space="preserve">part, </w:t></w:r><w:r><w:t xml:space="preserve">since </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">for </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>loops</w:t></w:r><w:r><w:t xml:space="preserve"> can be a little confusing, but learning them is worth it as they are incredibly useful.  However</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">there are other more skillful methods that avoid </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>for loops</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>For example,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>apply</w:t></w:r><w:r><w:t xml:space="preserve"> could replace </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> loop</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(16).Range.InsertXML($pkgHeader + $p16 + $pkgFooter)

$p17 = @'
<w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>emember:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>If I have a growth rate</w:t></w:r><w:r><w:t xml:space="preserve"> and multiply it with the </w:t></w:r><w:r><w:t xml:space="preserve">population </w:t></w:r><w:r><w:t xml:space="preserve">size </w:t></w:r><w:r><w:t>of a year</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I get the population </w:t></w:r><w:r><w:t xml:space="preserve">size </w:t></w:r><w:r><w:t xml:space="preserve">of the </w:t></w:r><w:r><w:t>following</w:t></w:r><w:r><w:t xml:space="preserve"> year</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">You </w:t></w:r><w:r><w:t xml:space="preserve">can use that knowledge to calculate the growth rate if you have the population </w:t></w:r><w:r><w:t xml:space="preserve">size </w:t></w:r><w:r><w:t xml:space="preserve">of </w:t></w:r><w:r><w:t xml:space="preserve">any </w:t></w:r><w:r><w:t xml:space="preserve">two </w:t></w:r><w:r><w:t xml:space="preserve">consecutive </w:t></w:r><w:r><w:t xml:space="preserve">years. </w:t></w:r></w:p>
'@
$d.Paragraphs(17).Range.InsertXML($pkgHeader + $p17 + $pkgFooter)

$p19 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Possible function</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>min(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>), max(</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>)</w:t></w:r></w:p>
'@
$d.Paragraphs(19).Range.InsertXML($pkgHeader + $p19 + $pkgFooter)

$p20 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Note:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>As any mathematical function in r (</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>mean</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>sum</w:t></w:r><w:r><w:t xml:space="preserve">…) if you have any NA’s in your data (only the column you calling in the function) </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">min </w:t></w:r><w:r><w:t xml:space="preserve">or </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>max</w:t></w:r><w:r><w:t xml:space="preserve"> will return NA</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$d.Paragraphs(20).Range.InsertXML($pkgHeader + $p20 + $pkgFooter)

$p23 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Possible function: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>prod(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…)</w:t></w:r></w:p>
'@
$d.Paragraphs(23).Range.InsertXML($pkgHeader + $p23 + $pkgFooter)

$p28 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Possible function:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>hist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…)</w:t></w:r></w:p>
'@
$d.Paragraphs(28).Range.InsertXML($pkgHeader + $p28 + $pkgFooter)

$p32 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Possible functions</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>sample(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>…), for(…)</w:t></w:r></w:p>
'@
$d.Paragraphs(32).Range.InsertXML($pkgHeader + $p32 + $pkgFooter)

$p35 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Plot your future population size t</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>h</w:t></w:r><w:r><w:t xml:space="preserve">rough time. </w:t></w:r></w:p>
'@
$d.Paragraphs(35).Range.InsertXML($pkgHeader + $p35 + $pkgFooter)

$p44 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Turn in the histogram of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>6</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">.) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>9</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>.)</w:t></w:r></w:p>
'@
$d.Paragraphs(44).Range.InsertXML($pkgHeader + $p44 + $pkgFooter)

$p38 = '<w:p/>'
$d.Paragraphs(38).Range.InsertXML($pkgHeader + $p38 + $pkgFooter)
